$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.520.21'
Set-TextValue 'E2' '  +2.22%  '
Set-TextValue 'D3' '1.859.39'
Set-TextValue 'E3' '  +1.55%  '
Set-TextValue 'D4' '0.9992'
Set-TextValue 'E4' '  -0.08%  '
Set-TextValue 'D5' '245.73'
Set-TextValue 'E5' '  +0.31%  '
Set-TextValue 'D6' '0.6964'
Set-TextValue 'E6' '  +1.11%  '
Set-TextValue 'D7' '0.9997'
Set-TextValue 'E7' '  -0.08%  '
Set-TextValue 'D8' '0.3080'
Set-TextValue 'E8' '  +1.04%  '
Set-TextValue 'D9' '0.07704'
Set-TextValue 'E9' '  +0.32%  '
Set-TextValue 'D10' '23.68'
Set-TextValue 'E10' '  +1.40%  '
Set-TextValue 'E11' '  -0.16%  '
Set-TextValue 'D12' '5.166'
Set-TextValue 'E12' '  +1.65%  '
Set-TextValue 'D13' '1.832.75'
Set-TextValue 'E13' '  +0.04%  '
Set-TextValue 'D14' '0.6958'
Set-TextValue 'E14' '  +2.12%  '
Set-TextValue 'D15' '91.23'
Set-TextValue 'E15' '  +1.08%  '
Set-TextValue 'D16' '6.361'
Set-TextValue 'E16' '  -1.26%  '
Set-TextValue 'D17' '29.492.60'
Set-TextValue 'E17' '  +2.17%  '
Set-TextValue 'D18' '0.000008320'
Set-TextValue 'E18' '  +0.27%  '
Set-TextValue 'D19' '2.100.25'
Set-TextValue 'E19' '  +1.36%  '
Set-TextValue 'D20' '238.42'
Set-TextValue 'E20' '  -1.31%  '
Set-TextValue 'D21' '12.77'
Set-TextValue 'E21' '  +0.55%  '
Set-TextValue 'D22' '0.9994'
Set-TextValue 'E22' '  -0.07%  '
Set-TextValue 'D23' '7.634'
Set-TextValue 'E23' '  +2.38%  '
Set-TextValue 'D24' '0.9998'
Set-TextValue 'E24' '  -0.15%  '
Set-TextValue 'E25' '  +1.40%  '
Set-TextValue 'B26' 'Cosmos'
Set-TextValue 'C26' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D26' '8.910'
Set-TextValue 'E26' '  +1.28%  '
Set-TextValue 'B27' 'Monero'
Set-TextValue 'C27' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D27' '160.05'
Set-TextValue 'E27' '  -0.73%  '
Set-TextValue 'D28' '18.29'
Set-TextValue 'E28' '  +0.64%  '
Set-TextValue 'E29' '  -0.78%  '
Set-TextValue 'D30' '4.254'
Set-TextValue 'E30' '  +1.06%  '
Set-TextValue 'D31' '4.155'
Set-TextValue 'E31' '  +0.12%  '
Set-TextValue 'E32' '  +2.11%  '
Set-TextValue 'D33' '0.05112'
Set-TextValue 'E33' '  +0.17%  '
Set-TextValue 'D34' '0.7784'
Set-TextValue 'E34' '  +1.62%  '
Set-TextValue 'D35' '1.887'
Set-TextValue 'E35' '  +2.76%  '
Set-TextValue 'D36' '1.150'
Set-TextValue 'E36' '  +1.03%  '
Set-TextValue 'D37' '2.687'
Set-TextValue 'E37' '  -0.28%  '
Set-TextValue 'D38' '1.318.08'
Set-TextValue 'E38' '  +7.96%  '
Set-TextValue 'D39' '0.01876'
Set-TextValue 'E39' '  +1.62%  '
Set-TextValue 'D40' '2.727'
Set-TextValue 'E40' '  +1.11%  '
Set-TextValue 'D41' '0.9518'
Set-TextValue 'E41' '  +1.18%  '
Set-TextValue 'D42' '105.94'
Set-TextValue 'E42' '  -2.44%  '
Set-TextValue 'E43' '  +1.04%  '
Set-TextValue 'D44' '1.000'
Set-TextValue 'E44' '  +0.04%  '
Set-TextValue 'D45' '9.849'
Set-TextValue 'E45' '  +3.70%  '
Set-TextValue 'E46' '  +2.00%  '
Set-TextValue 'D47' '1.999.52'
Set-TextValue 'E47' '  +1.23%  '
Set-TextValue 'D48' '0.5231'
Set-TextValue 'E48' '  +1.34%  '
Set-TextValue 'D49' '1.792'
Set-TextValue 'E49' '  +2.71%  '
Set-TextValue 'D50' '63.12'
Set-TextValue 'E50' '  -1.35%  '
Set-TextValue 'D51' '6.976'
Set-TextValue 'E51' '  +1.37%  '
